# New crime data collected
# Applies the weekly cs-en-us-032pct.xlsx data refresh described by the
# commit's canonical OOXML diff: new volume/report-week header text, and
# refreshed counts / % change figures in rows 14-29 of the CompStat table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text runs inside shared strings) - only the numeric
# run segments change, surrounding label runs stay as-is.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  16" -> "...  17"
$a8 = $ws.Range("A8")
$a8full = $a8.Value2
$a8.Characters($a8full.Length - 1, 2).Text = "17"

# C9: "Report Covering the Week  4/17/2023  Through  4/23/2023"
#     -> "...4/24/2023  Through  4/30/2023"
$c9 = $ws.Range("C9")
$c9full = $c9.Value2
$firstDateStart = $c9full.IndexOf("4/17/2023") + 1
$c9.Characters($firstDateStart, 9).Text = "4/24/2023"
$c9full2 = $ws.Range("C9").Value2
$secondDateStart = $c9full2.IndexOf("4/23/2023") + 1
$ws.Range("C9").Characters($secondDateStart, 9).Text = "4/30/2023"

# ---------------------------------------------------------------------
# Data table refresh, rows 14-29
# ---------------------------------------------------------------------

# Row 14
$ws.Range("N14").Value2 = -76.923076923076

# Row 15
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 1
$ws.Range("H15").Value2 = 100
$ws.Range("M15").Value2 = 0
$ws.Range("N15").Value2 = -65.517241379310

# Row 16
$ws.Range("C16").Value2 = 7
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 600
$ws.Range("F16").Value2 = 13
$ws.Range("G16").Value2 = 21
$ws.Range("H16").Value2 = -38.095238095238
$ws.Range("I16").Value2 = 52
$ws.Range("J16").Value2 = 88
$ws.Range("K16").Value2 = -40.909090909090
$ws.Range("L16").Value2 = 36.842105263157
$ws.Range("M16").Value2 = -40.909090909090
$ws.Range("N16").Value2 = -81.690140845070

# Row 17
$ws.Range("C17").Value2 = 6
$ws.Range("D17").Value2 = 12
$ws.Range("E17").Value2 = -50
$ws.Range("F17").Value2 = 32
$ws.Range("G17").Value2 = 49
$ws.Range("H17").Value2 = -34.693877551020
$ws.Range("I17").Value2 = 141
$ws.Range("J17").Value2 = 169
$ws.Range("K17").Value2 = -16.568047337278
$ws.Range("L17").Value2 = 0
$ws.Range("M17").Value2 = 65.882352941176
$ws.Range("N17").Value2 = -50.526315789473

# Row 18
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 2
$ws.Range("E18").Value2 = 50
$ws.Range("F18").Value2 = 9
$ws.Range("G18").Value2 = 19
$ws.Range("H18").Value2 = -52.631578947368
$ws.Range("I18").Value2 = 57
$ws.Range("J18").Value2 = 82
$ws.Range("K18").Value2 = -30.487804878048
$ws.Range("L18").Value2 = 1.785714285714
$ws.Range("M18").Value2 = 62.857142857142
$ws.Range("N18").Value2 = -81.904761904761

# Row 19
$ws.Range("C19").Value2 = 3
$ws.Range("D19").Value2 = 1
$ws.Range("E19").Value2 = 200
$ws.Range("G19").Value2 = 23
$ws.Range("H19").Value2 = 34.782608695652
$ws.Range("I19").Value2 = 106
$ws.Range("J19").Value2 = 106
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 11.578947368421
$ws.Range("M19").Value2 = 16.483516483516
$ws.Range("N19").Value2 = -13.821138211382

# Row 20 - C20 switches from a numeric 4 to a literal text "0" (quote
# prefix), matching the other zero-count cells elsewhere in the sheet
# (e.g. C14/D14). Re-apply font + alignment so it still reads like the
# surrounding number cells.
$ws.Range("C20").Value2 = "'0"
$ws.Range("C20").Font.Name = "Andale WT"
$ws.Range("C20").Font.Size = 10
$ws.Range("C20").HorizontalAlignment = -4152
$ws.Range("C20").VerticalAlignment = -4108
$ws.Range("C20").NumberFormat = "General"
$ws.Range("E20").Value2 = -100
$ws.Range("F20").Value2 = 10
$ws.Range("G20").Value2 = 10
$ws.Range("H20").Value2 = 0
$ws.Range("J20").Value2 = 31
$ws.Range("K20").Value2 = -12.903225806451
$ws.Range("M20").Value2 = 80
$ws.Range("N20").Value2 = -72.164948453608

# Row 21 (TOTAL)
$ws.Range("C21").Value2 = 19
$ws.Range("D21").Value2 = 18
$ws.Range("E21").Value2 = 5.555555555555
$ws.Range("F21").Value2 = 98
$ws.Range("G21").Value2 = 123
$ws.Range("H21").Value2 = -20.325203252032
$ws.Range("I21").Value2 = 396
$ws.Range("J21").Value2 = 488
$ws.Range("K21").Value2 = -18.852459016393
$ws.Range("L21").Value2 = 12.5
$ws.Range("M21").Value2 = 21.846153846153
$ws.Range("N21").Value2 = -65.445026178010

# Row 22
$ws.Range("G22").Value2 = 2
$ws.Range("H22").Value2 = -50

# Row 23
$ws.Range("C23").Value2 = 4
$ws.Range("E23").Value2 = 100
$ws.Range("F23").Value2 = 17
$ws.Range("G23").Value2 = 12
$ws.Range("H23").Value2 = 41.666666666666
$ws.Range("I23").Value2 = 70
$ws.Range("J23").Value2 = 67
$ws.Range("K23").Value2 = 4.477611940298
$ws.Range("L23").Value2 = 11.111111111111
$ws.Range("M23").Value2 = 40

# Row 24
$ws.Range("C24").Value2 = 7
$ws.Range("D24").Value2 = 18
$ws.Range("E24").Value2 = -61.111111111111
$ws.Range("F24").Value2 = 54
$ws.Range("H24").Value2 = -18.181818181818
$ws.Range("I24").Value2 = 262
$ws.Range("J24").Value2 = 283
$ws.Range("K24").Value2 = -7.420494699646
$ws.Range("L24").Value2 = 11.965811965812
$ws.Range("M24").Value2 = 25.358851674641

# Row 25
$ws.Range("C25").Value2 = 22
$ws.Range("D25").Value2 = 11
$ws.Range("E25").Value2 = 100
$ws.Range("F25").Value2 = 51
$ws.Range("G25").Value2 = 56
$ws.Range("H25").Value2 = -8.928571428571
$ws.Range("I25").Value2 = 192
$ws.Range("J25").Value2 = 193
$ws.Range("K25").Value2 = -0.518134715025
$ws.Range("L25").Value2 = 24.675324675324
$ws.Range("M25").Value2 = -27.272727272727

# Row 26 - D26/E26 switch from literal text ("0" / "***.*") to real
# numbers. Re-apply the usual numeric-cell number formats so the style
# lines up with the rest of the column (D -> #,##0 ; E -> the signed
# one-decimal % format).
$ws.Range("D26").Value2 = 2
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value2 = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F26").Value2 = 3
$ws.Range("G26").Value2 = 3
$ws.Range("H26").Value2 = 0
$ws.Range("J26").Value2 = 13
$ws.Range("K26").Value2 = 15.384615384615

# Row 27
$ws.Range("C27").Value2 = 2
$ws.Range("D27").Value2 = 4
$ws.Range("E27").Value2 = -50
$ws.Range("F27").Value2 = 9
$ws.Range("G27").Value2 = 11
$ws.Range("H27").Value2 = -18.181818181818
$ws.Range("I27").Value2 = 21
$ws.Range("J27").Value2 = 27
$ws.Range("K27").Value2 = -22.222222222222
$ws.Range("L27").Value2 = -19.230769230769

# Row 28
$ws.Range("J28").Value2 = 10
$ws.Range("K28").Value2 = -20
$ws.Range("N28").Value2 = -76.470588235294

# Row 29
$ws.Range("J29").Value2 = 9
$ws.Range("K29").Value2 = -22.222222222222
$ws.Range("N29").Value2 = -76.666666666666
